$wb = $excel.ActiveWorkbook

# The order/source file dropped into row 3 ("fileName" column A) was
# re-uploaded as a newer revision: stp-testing 3.stp -> stp-testing 4.stp.
# This same row 3 "fileName" value is duplicated across every order sheet.
$sheetNames = @(
    "FPA001",
    "FPA002-003-005-007",
    "FPA004-006-010",
    "FPA008-009",
    "BTMI002",
    "BTMI003",
    "BTMI015"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A3").Value = "stp-testing 4.stp"
}

# Re-create the click-through selection trail left behind while reviewing
# each sheet's row 3, finishing back on FPA001 so it ends up the active tab.
$wb.Worksheets.Item("FPA002-003-005-007").Range("A3").Select()
$wb.Worksheets.Item("FPA004-006-010").Range("A3").Select()
$wb.Worksheets.Item("FPA008-009").Range("A3").Select()
$wb.Worksheets.Item("BTMI002").Range("A3").Select()
$wb.Worksheets.Item("BTMI003").Range("A3").Select()
$wb.Worksheets.Item("BTMI015").Range("F26").Select()
$wb.Worksheets.Item("FPA001").Range("A3").Select()
